# =====================================================================
# 杭州-漫展信息.xlsx -- apply "456a3b4" gh-pages data refresh
# =====================================================================
$wb = $excel.ActiveWorkbook

# ---- Sheet "展览": refresh "想去人数" (column F) visit counters ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(2,6).Value = 258
$ws1.Cells.Item(3,6).Value = 2779
$ws1.Cells.Item(5,6).Value = 952
$ws1.Cells.Item(6,6).Value = 39
$ws1.Cells.Item(7,6).Value = 2960
$ws1.Cells.Item(8,6).Value = 1889
$ws1.Cells.Item(9,6).Value = 237
$ws1.Cells.Item(11,6).Value = 2543
$ws1.Cells.Item(12,6).Value = 574
$ws1.Cells.Item(13,6).Value = 263
$ws1.Cells.Item(16,6).Value = 139
$ws1.Cells.Item(18,6).Value = 9480
$ws1.Cells.Item(20,6).Value = 4
$ws1.Cells.Item(21,6).Value = 7462
$ws1.Cells.Item(22,6).Value = 11995
$ws1.Cells.Item(25,6).Value = 245
$ws1.Cells.Item(27,6).Value = 578
$ws1.Cells.Item(28,6).Value = 2698
$ws1.Cells.Item(29,6).Value = 241
$ws1.Cells.Item(30,6).Value = 213
$ws1.Cells.Item(31,6).Value = 2682
$ws1.Cells.Item(32,6).Value = 986
$ws1.Cells.Item(33,6).Value = 5
$ws1.Cells.Item(36,6).Value = 4549
$ws1.Cells.Item(37,6).Value = 1078
$ws1.Cells.Item(38,6).Value = 30
$ws1.Cells.Item(40,6).Value = 60
$ws1.Cells.Item(41,6).Value = 563

# ---- Sheet "本地生活": refresh F4 counter ----
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(4,6).Value = 188

# ---- Sheet "全部类型": refresh "想去人数" (column F) visit counters ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(3,6).Value = 258
$ws4.Cells.Item(5,6).Value = 2779
$ws4.Cells.Item(8,6).Value = 952
$ws4.Cells.Item(9,6).Value = 39
$ws4.Cells.Item(11,6).Value = 2960
$ws4.Cells.Item(13,6).Value = 1889
$ws4.Cells.Item(14,6).Value = 237
$ws4.Cells.Item(15,6).Value = 2543
$ws4.Cells.Item(17,6).Value = 574
$ws4.Cells.Item(18,6).Value = 263
$ws4.Cells.Item(20,6).Value = 139
$ws4.Cells.Item(22,6).Value = 9481
$ws4.Cells.Item(24,6).Value = 4
$ws4.Cells.Item(25,6).Value = 7462
$ws4.Cells.Item(26,6).Value = 11995
$ws4.Cells.Item(29,6).Value = 245
$ws4.Cells.Item(32,6).Value = 578
$ws4.Cells.Item(34,6).Value = 2698
$ws4.Cells.Item(36,6).Value = 241
$ws4.Cells.Item(37,6).Value = 213
$ws4.Cells.Item(40,6).Value = 4549
$ws4.Cells.Item(45,6).Value = 563

# ---- Sheet "演出": new concert announced for 2024-12-27 ----
$ws2 = $wb.Worksheets.Item("演出")

# existing row 22 "想去人数" ticked up independently of the insert below
$ws2.Cells.Item(22,6).Value = 18

# make room at row 23; rows 23-27 (the events that follow chronologically)
# shift down to 24-28
$ws2.Rows.Item(23).Insert()

# running index column keeps counting 0,1,2,... down column A regardless of
# which event occupies the row
$ws2.Cells.Item(23,1).Value = 22

# column B holds a literal "yyyy-MM-dd" label, not a real date -- force text
# so Excel does not reinterpret it as a date serial
$bCell = $ws2.Cells.Item(23,2)
$bCell.NumberFormat = "@"
$bCell.Value = "2024-12-27"

$ws2.Cells.Item(23,3).Value = "杭州·德国美因茨名家管弦乐团2025新年音乐会"
$ws2.Cells.Item(23,4).Value = "建国南路280号（城站火车站附近） 杭州红星剧院"
$ws2.Cells.Item(23,5).Value = "2024.12.27 19:30-12.27 21:00"
$ws2.Cells.Item(23,6).Value = 0
$ws2.Cells.Item(23,7).Value = 126
$ws2.Cells.Item(23,8).Value = "https://show.bilibili.com/platform/detail.html?id=93583"
$ws2.Cells.Item(23,9).Value = "//i2.hdslb.com/bfs/openplatform/202410/EdEriHh21729143553031.jpeg"

# match the bold / centred / bordered look used by every other cell in the
# index column
$a23 = $ws2.Cells.Item(23,1)
$a23.Font.Bold = $true
$a23.HorizontalAlignment = -4108
$a23.VerticalAlignment = -4160
$a23.Borders.LineStyle = 1

# Insert() drags the old index numbers down along with rows 23-27 (now
# 24-28); column A is really just "row number - 1" everywhere in this
# sheet, so restate it for every row the insert displaced
for ($r = 24; $r -le 28; $r++) {
    $ws2.Cells.Item($r,1).Value = $r - 1
}
